$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (count)
$ws.Range("B2").Value = 50000
$ws.Range("C2").Value = 50000
$ws.Range("D2").Value = 50000
$ws.Range("E2").Value = 50000
$ws.Range("F2").Value = 50000
$ws.Range("G2").Value = 50000

# Row 3 (mean)
$ws.Range("B3").Value = 100.0000012
$ws.Range("C3").Value = 0.0999999964
$ws.Range("D3").Value = 0.0009999999280000001
$ws.Range("E3").Value = 0.7101831633425429
$ws.Range("F3").Value = 1.034148332520787
$ws.Range("G3").Value = 80582.32065237426

# Row 4 (std)
$ws.Range("B4").Value = 17.33801102677369
$ws.Range("C4").Value = 0.01733799088564029
$ws.Range("D4").Value = 0.0001733800567587762
$ws.Range("E4").Value = 0.04184318341006341
$ws.Range("F4").Value = 0.003428288513016337
$ws.Range("G4").Value = 13328.3073714289

# Row 5 (min)
$ws.Range("E5").Value = 0.5721519384014822
$ws.Range("F5").Value = 1.029535005472516
$ws.Range("G5").Value = 55230.42756676323

# Row 6 (25%)
$ws.Range("E6").Value = 0.6817977762586498
$ws.Range("F6").Value = 1.031581962698854
$ws.Range("G6").Value = 69782.59832743558

# Row 7 (50%)
$ws.Range("E7").Value = 0.7102005413632267
$ws.Range("F7").Value = 1.033275864086499
$ws.Range("G7").Value = 78769.4728925252

# Row 8 (75%)
$ws.Range("E8").Value = 0.7408790013968939
$ws.Range("F8").Value = 1.035737458526838
$ws.Range("G8").Value = 90054.19338616684

# Row 9 (max)
$ws.Range("E9").Value = 0.8108401949089699
$ws.Range("F9").Value = 1.049419532824653
$ws.Range("G9").Value = 126839.118783051
